$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "90.757.90"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.099.02"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.89"
$ws.Range("E5").Value = "  +4.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "620.46"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.15"
$ws.Range("E7").Value = "  +4.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.365"
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.099.74"
$ws.Range("E10").Value = "  +24.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.742"
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("E12").Value = "  +3.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.18"
$ws.Range("E14").Value = "  -3.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.48"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.383.55"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.669.55"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.101.74"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.78"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.66"
$ws.Range("E20").Value = "  +4.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000212"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.82"
$ws.Range("E22").Value = "  +4.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "440.48"
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.04"
$ws.Range("E24").Value = "  +1.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.67"
$ws.Range("E25").Value = "  -3.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "90.05"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.90"
$ws.Range("E27").Value = "  -2.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.245.05"
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.247"
$ws.Range("E30").Value = "  +22.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.177"
$ws.Range("E31").Value = "  +11.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.121"
$ws.Range("E32").Value = "  +36.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.25"
$ws.Range("E33").Value = "  -2.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.168"
$ws.Range("E34").Value = "  +9.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.966"
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.94"
$ws.Range("E36").Value = "  +13.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.34"
$ws.Range("E37").Value = "  +24.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.30"
$ws.Range("E39").Value = "  +0.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "493.41"
$ws.Range("E40").Value = "  -2.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.60"
$ws.Range("E41").Value = "  -6.02%  "
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("E43").Value = "  +2.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.17"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.92"
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "155.17"
$ws.Range("E47").Value = "  +2.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.689"
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.24"
$ws.Range("E50").Value = "  -1.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.38"
$ws.Range("E51").Value = "  -1.27%  "
